$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.812.48'

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.21%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.888.57'

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.02%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '566.82'

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -3.76%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.93'

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.56%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.10%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.501'

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.80%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.885.44'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.93'

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.14%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.147'

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.11%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.08%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.00%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '31.85'

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.61%  '

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.04%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.367.86'

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.94%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.774.47'

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.27%  '

$ws.Range("B18").Value = 'WrappedEther'

$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.894.99'

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.78%  '

$ws.Range("B19").Value = 'Polkadot'

$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.52'

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.57%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '429.09'

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.27%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.99'

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -3.02%  '

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.91%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.85'

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.22%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.91'

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.37%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.01'

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.69%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.04'

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -9.97%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.10%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.30%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000109'

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +9.63%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.96'

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.75%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.15%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.01'

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -7.07%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.84%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.06%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.94%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.949'

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.94%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.37'

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.07%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '48.83'

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.65%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.80'

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -6.11%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.38%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.115'

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.15%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.15'

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.27%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '39.68'

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.41%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.268'

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.99%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.687.45'

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.42%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.90%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '131.42'

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.35%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '344.14'

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.08%  '

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.03%  '

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.25%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '21.51'

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.91%  '
